$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Port_Letter")

# Delete the now-unused trailing empty rows of the header table (rows 42-44 first,
# bottom-up, so row numbers stay stable while we work), then the extra blank rows
# right after the header table (rows 29-32).
$ws.Rows("42:44").Delete()
$ws.Rows("29:32").Delete()

# The two conditional-formatting rules still cover their old (now too-large)
# ranges after the row deletes; shrink them to match the shrunk table.
$ws.Range("A24:F31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("A24:F28"))
$ws.Range("H29:N48").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H29:N41"))

# Update the print scale
$ws.PageSetup.Zoom = 86

# Restore view state (selection / scroll position) to match the saved file
$ws.Range("A32:F32").Select()
$ws.Application.ActiveWindow.ScrollRow = 20
